$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 164, pushing existing rows 164-197 down to 165-198.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new data entry.
$ws.Cells.Item(164, 1).Value  = 11
$ws.Cells.Item(164, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(164, 3).Value  = "Bíobío"
$ws.Cells.Item(164, 4).Value  = 44798
$ws.Cells.Item(164, 5).Value  = 8
$ws.Cells.Item(164, 6).Value  = "Fruta"
$ws.Cells.Item(164, 7).Value  = 100108
$ws.Cells.Item(164, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value  = 100108005
$ws.Cells.Item(164, 10).Value = "Piña"
$ws.Cells.Item(164, 11).Value = "Caramelo"
$ws.Cells.Item(164, 12).Value = "Segunda"
$ws.Cells.Item(164, 13).Value = 220
$ws.Cells.Item(164, 14).Value = 18000
$ws.Cells.Item(164, 15).Value = 19000
$ws.Cells.Item(164, 16).Value = 18455
$ws.Cells.Item(164, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(164, 18).Value = "Ecuador"
$ws.Cells.Item(164, 19).Value = 1318
$ws.Cells.Item(164, 20).Value = 14
